# Third commit for 4th march 2017
# - "day 1" sheet: finish filling in row 69 (S.No 68) and add new row 70 (S.No 69)
# - "error report" sheet: add a new error/solution entry in row 10

$wb = $excel.ActiveWorkbook

$day1 = $wb.Worksheets.Item("day 1")
$errorReport = $wb.Worksheets.Item("error report")

# ---- "day 1" sheet ----

# Row 69 (S.No 68) - "Write test cases to work with ProductDAO" already has A/B/C filled in;
# finish the rest of the row.
$day1.Range("D69").Value = "https://www.youtube.com/watch?v=S9wKr2nuHHw&t=25s"
$day1.Range("E69").Value = "NA"
$day1.Range("F69").Value = "60 minutes"
$day1.Range("G69").Value = "N"
$day1.Range("H69").Value = "NA"

# Row 70 (S.No 69) - new "Documentation" task row.
$day1.Range("B70").Value = "4th Mar,2017"
$day1.Range("C70").Value = "Documentation"
$day1.Range("D70").Value = "NA"
$day1.Range("E70").Value = "NA"
$day1.Range("F70").Value = "30 minutes"
$day1.Range("G70").Value = "N"
$day1.Range("H70").Value = "NA"
$day1.Rows.Item(70).RowHeight = 28.8

# ---- "error report" sheet ----

# Row 10 - new error report entry about images not loading.
$errorReport.Range("A10").Value = "Images_Error"
$errorReport.Range("A10").WrapText = $true
$errorReport.Range("B10").Value = "Images were not displayed after doing all the static resource loading"
$errorReport.Range("C10").Value = "Ctrl+F5 in the browser , this clears the cache"
$errorReport.Range("C10").WrapText = $true
$errorReport.Range("D10").Value = "NA"
$errorReport.Range("D10").WrapText = $true
$errorReport.Rows.Item(10).RowHeight = 28.8

# Restore the final selection / scroll position on the "day 1" sheet.
[void]$day1.Activate()
[void]$day1.Range("H70").Select()
